$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20, shifting existing rows 20-65
# down to 21-66 (matches the diff: old row 20 data now lives at row 21, ...,
# old row 65 data now lives at row 66).
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new weekly data point.
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C20").Value = "Metropolitana"
$ws.Range("D20").Value = 44715
$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 100112029
$ws.Range("G20").Value = "Orégano"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 14000
$ws.Range("N20").Value = "$/docena de atados"
$ws.Range("O20").Value = "Región Metropolitana"
$ws.Range("P20").Value = 4667
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Hortaliza"
